$wb = $excel.ActiveWorkbook

# Add the new worksheet at the very end of the sheet collection.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "zoom150"

# Populate A1 with the descriptive text.
$newSheet.Range("A1").Value = "Zoom scale 150%"

# Set the zoom level for this sheet's view to 150%.
$newSheet.Activate()
$excel.ActiveWindow.Zoom = 150

# Adding/naming a new sheet shifts the workbook's "active tab" to it; put
# the active tab back on the "active" sheet (workbook's activeTab=5, i.e.
# the 6th sheet) so this edit only adds the new sheet/data, leaving the
# previously-active tab selection untouched.
$wb.Worksheets.Item("active").Activate()
